# Ticket 46 - add a new "jt:rickroll" example row to the Hyperlinks sheet
# of the JETT HyperlinkTagTemplate workbook. This appends a new shared
# string (rich text, 3 runs) in cell A3:
#   <jt:rickroll value="   -- plain text
#   Additional Help         -- blue, underlined ("hyperlink" look)
#   "/>                     -- plain text
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$prefix = '<jt:rickroll value="'
$label  = 'Additional Help'
$suffix = '"/>'
$full   = $prefix + $label + $suffix

$cell = $ws.Range("A3")
$cell.Value = $full

# Style the "Additional Help" run like a hyperlink (blue + underline),
# matching the formatting used for the other label runs on this sheet.
$labelStart = $prefix.Length + 1
$labelRun = $cell.Characters($labelStart, $label.Length)
$labelRun.Font.Underline = $true
$labelRun.Font.Color = 16711680

# Give the trailing `"/>` run its own explicit (default-looking) font
# properties, matching the existing rich-text cells on this sheet.
$suffixStart = $prefix.Length + $label.Length + 1
$suffixRun = $cell.Characters($suffixStart, $suffix.Length)
$suffixRun.Font.Name = "Calibri"
$suffixRun.Font.Size = 11
